# Regenerate the "K" column (G) values for the save_data sheet.
# The commit replaces the previous Strike# derived values in column G
# (rows 2-36) with freshly recalculated K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(4,7,2,5,4,9,2,6,5,6,7,4,7,3,11,7,2,6,5,2,5,7,6,4,1,7,8,3,6,2,7,5,1,2,3)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
